# Update gh-pages output data (new crawl snapshot at 456a3b4)
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10667
$ws1.Range("F3").Value = 237
$ws1.Range("F4").Value = 62
$ws1.Range("F5").Value = 690
$ws1.Range("F6").Value = 494

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10667
$ws4.Range("F3").Value = 237
$ws4.Range("F4").Value = 62
$ws4.Range("F5").Value = 690
$ws4.Range("F7").Value = 494
